# poisson_naive pronta para a rodada 27
# Updates the "A" (game id) and "E" (Round) columns of the Gremio fixtures
# sheet, and rotates the full data of rows 22-24 by one row (22->23->24->22)
# as described by the source diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: snapshot the full row contents (columns B..BD) for rows 22,23,24
# before we start writing anything, since row 22's new data comes from row 24's
# old data, row 23's new data comes from row 22's old data, and row 24's new
# data comes from row 23's old data (a downward rotation).

$lastCol = 58  # column BD

function Get-RowValues($rowIndex) {
    $vals = @{}
    for ($c = 2; $c -le $lastCol; $c++) {
        $cell = $ws.Cells.Item($rowIndex, $c)
        $vals[$c] = $cell.Value()
    }
    return $vals
}

$row22 = Get-RowValues 22
$row23 = Get-RowValues 23
$row24 = Get-RowValues 24

function Set-RowValues($rowIndex, $vals) {
    for ($c = 2; $c -le $lastCol; $c++) {
        $cell = $ws.Cells.Item($rowIndex, $c)
        $val = $vals[$c]
        if ($c -eq 2) {
            # Column B holds dates stored as plain text (e.g. "2023-04-22").
            # Force a text number format first so Excel does not silently
            # convert the string into a date serial number.
            $cell.NumberFormat = "@"
        }
        if ($null -eq $val) {
            $cell.Value = ""
        } else {
            $cell.Value = $val
        }
    }
}

# row22 <- old row24 ; row23 <- old row22 ; row24 <- old row23
Set-RowValues 22 $row24
Set-RowValues 23 $row22
Set-RowValues 24 $row23

# --- Step 2: set the final A (col 1) and E (col 5) values for every data row.
# E must end up as a plain number (it used to be text like "Matchweek N").

$targets = @{
    2  = @{ A = 0;  E = 1 }
    3  = @{ A = 3;  E = 4 }
    4  = @{ A = 5;  E = 6 }
    5  = @{ A = 6;  E = 7 }
    6  = @{ A = 8;  E = 9 }
    7  = @{ A = 10; E = 11 }
    8  = @{ A = 11; E = 12 }
    9  = @{ A = 13; E = 14 }
    10 = @{ A = 14; E = 16 }
    11 = @{ A = 17; E = 19 }
    12 = @{ A = 19; E = 21 }
    13 = @{ A = 20; E = 22 }
    14 = @{ A = 23; E = 24 }
    15 = @{ A = 30; E = 23 }
    16 = @{ A = 7;  E = 5 }
    17 = @{ A = 14; E = 10 }
    18 = @{ A = 39; E = 25 }
    19 = @{ A = 11; E = 8 }
    20 = @{ A = 2;  E = 3 }
    21 = @{ A = 37; E = 26 }
    22 = @{ A = 1;  E = 2 }
    23 = @{ A = 34; E = 15 }
    24 = @{ A = 25; E = 20 }
    25 = @{ A = 12; E = 13 }
    26 = @{ A = 16; E = 18 }
    27 = @{ A = 22; E = 17 }
}

foreach ($r in $targets.Keys) {
    $t = $targets[$r]
    $ws.Cells.Item($r, 1).Value = $t.A
    $eCell = $ws.Cells.Item($r, 5)
    $eCell.NumberFormat = "General"
    $eCell.Value = $t.E
}
